$d = $word.ActiveDocument

# Locate the paragraph that holds the three "Requisitos" list entries by
# searching for the LOB1012 requirement line (it is currently the first
# run of that paragraph).
$search = $d.Content
$search.Find.ClearFormatting()
$found = $search.Find.Execute("LOB1012 -  Estatística  (Requisito fraco)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the LOB1012 requirement line"
}

$targetParagraph = $search.Paragraphs(1)
$pStart = $targetParagraph.Range.Start
$pEnd = $targetParagraph.Range.End

# Range covering the paragraph's runs but excluding the trailing paragraph
# mark, so re-inserting OOXML here replaces only the run content and keeps
# the existing paragraph (with its ListBullet style) intact.
$runsRange = $d.Range($pStart, $pEnd - 1)

$newRunsXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>LOQ4095 -  Química Geral Experimental  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1012 -  Estatística  (Requisito fraco)</w:t><w:br/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$runsRange.InsertXML($newRunsXml)
